$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "time_taken" column, matching the bold/bordered
# header style already used by the other header cells in row 1 (B1:E1).
# Copy+PasteSpecial(Formats) reuses the existing cell style (s="1") rather
# than cloning a brand new style entry into styles.xml.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Timestamps (one per data row, 2-18) recorded as plain text values.
$timestamps = @(
    "2021-10-05 10:52:11.213037",
    "2021-10-05 10:52:11.213050",
    "2021-10-05 10:52:11.213053",
    "2021-10-05 10:52:11.213057",
    "2021-10-05 10:52:11.213060",
    "2021-10-05 10:52:11.213063",
    "2021-10-05 10:52:11.213066",
    "2021-10-05 10:52:11.213069",
    "2021-10-05 10:52:11.213073",
    "2021-10-05 10:52:11.213076",
    "2021-10-05 10:52:11.213079",
    "2021-10-05 10:52:11.213082",
    "2021-10-05 10:52:11.213085",
    "2021-10-05 10:52:11.213088",
    "2021-10-05 10:52:11.213091",
    "2021-10-05 10:52:11.213094",
    "2021-10-05 10:52:11.213097"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
